$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 205, shifting the existing rows 205:240 down to 206:241
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new data point
$ws.Cells.Item(205, 1).Value = 10
$ws.Cells.Item(205, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(205, 3).Value = "La Araucanía"
$ws.Cells.Item(205, 4).Value = 44951
$ws.Cells.Item(205, 5).Value = 9
$ws.Cells.Item(205, 6).Value = 100114007
$ws.Cells.Item(205, 7).Value = "Jengibre"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 25
$ws.Cells.Item(205, 11).Value = 22000
$ws.Cells.Item(205, 12).Value = 22000
$ws.Cells.Item(205, 13).Value = 22000
$ws.Cells.Item(205, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(205, 15).Value = "Perú"
$ws.Cells.Item(205, 16).Value = 1692
$ws.Cells.Item(205, 17).Value = 13
$ws.Cells.Item(205, 18).Value = "Hortaliza"
